# Fixed Stimulus Absolute Timestamps
$wb = $excel.ActiveWorkbook

# Rename sheets (task order sheet tabs) - refresh the embedded timestamps
$wb.Worksheets.Item(1).Name = "GNG_TO-16504778531284919"
$wb.Worksheets.Item(2).Name = "NB_TO-16504778552344966"
$wb.Worksheets.Item(3).Name = "RS_TO-1650477855235495"
$wb.Worksheets.Item(4).Name = "TOL_TO-16504778552974916"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16504778553614914"

# Sheet 1: GNG
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16504778530862436.csv"
$ws1.Range("B3").Value = "GNG_stims-16504778531116388.csv"
$ws1.Range("B4").Value = "go_stims-16504778531124935.csv"
$ws1.Range("B5").Value = "GNG_stims-1650477853127491.csv"

# Sheet 2: NB
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16504778541955273.csv"
$ws2.Range("B3").Value = "ZB-match_0-16504778533244958.csv"
$ws2.Range("B4").Value = "TB-1650477855120529.csv"
$ws2.Range("B5").Value = "OB-16504778542295265.csv"
$ws2.Range("B6").Value = "TB-1650477855217526.csv"
$ws2.Range("B7").Value = "OB-16504778544365296.csv"
$ws2.Range("B8").Value = "TB-16504778549625127.csv"
$ws2.Range("B9").Value = "ZB-match_3-1650477853554501.csv"
$ws2.Range("B10").Value = "ZB-match_5-16504778534564962.csv"

# Sheet 3: RS
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# Sheet 4: TOL
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16504778552654915.csv"
$ws4.Range("B3").Value = "ZM_stims-16504778552424932.csv"
$ws4.Range("B4").Value = "MM_stims-16504778552814913.csv"
$ws4.Range("B5").Value = "ZM_stims-16504778552664914.csv"
$ws4.Range("B6").Value = "MM_stims-16504778552974916.csv"
$ws4.Range("B7").Value = "ZM_stims-16504778552824926.csv"

# Sheet 5: vSAT
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16504778553294911.csv"
$ws5.Range("B3").Value = "SAT_stims-16504778553134913.csv"
$ws5.Range("B4").Value = "SAT_stims-16504778553004928.csv"
$ws5.Range("B5").Value = "vSAT_stims-16504778553454921.csv"
